# Auto-generated Excel COM-interop script to update crypto price/volume data
# Commit: Updated cryptos list on Tue Oct 31 14:42:06 UTC 2023 with GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.287.37'
$ws.Range('E2').Value = '  -1.17%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.796.74'
$ws.Range('E3').Value = '  -1.26%  '

# Row 4
$ws.Range('E4').Value = '  +0.19%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.19'
$ws.Range('E5').Value = '  -1.01%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.593'
$ws.Range('E6').Value = '  +2.17%  '

# Row 7
$ws.Range('E7').Value = '  +0.17%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '36.03'
$ws.Range('E8').Value = '  +3.83%  '

# Row 9
$ws.Range('E9').Value = '  -2.79%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0686'
$ws.Range('E10').Value = '  -2.56%  '

# Row 11
$ws.Range('E11').Value = '  +1.08%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.058.46'
$ws.Range('E12').Value = '  -1.14%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.25'
$ws.Range('E13').Value = '  -1.36%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.808.21'
$ws.Range('E14').Value = '  -0.70%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.636'
$ws.Range('E15').Value = '  -1.40%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '34.296.55'
$ws.Range('E16').Value = '  -1.11%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.39'
$ws.Range('E17').Value = '  +1.20%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.41'
$ws.Range('E18').Value = '  +0.30%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '242.59'
$ws.Range('E19').Value = '  -1.73%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0781'
$ws.Range('E20').Value = '  -2.86%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.36'
$ws.Range('E21').Value = '  -1.89%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.13'
$ws.Range('E23').Value = '  -1.35%  '

# Row 24
$ws.Range('E24').Value = '  +6.15%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '170.64'
$ws.Range('E25').Value = '  -2.01%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.96'
$ws.Range('E26').Value = '  +6.00%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.24'
$ws.Range('E27').Value = '  +2.38%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.120'
$ws.Range('E28').Value = '  +0.91%  '

# Row 29
$ws.Range('E29').Value = '  +0.18%  '

# Row 30
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.24'
$ws.Range('E30').Value = '  -0.73%  '

# Row 31
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.92'
$ws.Range('E31').Value = '  -2.03%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.80'
$ws.Range('E32').Value = '  -1.35%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0517'
$ws.Range('E33').Value = '  -2.81%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.78'
$ws.Range('E34').Value = '  -3.52%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.372.63'
$ws.Range('E35').Value = '  -2.62%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.651'
$ws.Range('E36').Value = '  -4.61%  '

# Row 37
$ws.Range('E37').Value = '  -2.15%  '

# Row 38
$ws.Range('E38').Value = '  -11.28%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0186'
$ws.Range('E39').Value = '  -3.43%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '81.65'
$ws.Range('E40').Value = '  -3.70%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.42'
$ws.Range('E41').Value = '  +0.75%  '

# Row 42
$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.79'
$ws.Range('E42').Value = '  -3.08%  '

# Row 43
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.940'
$ws.Range('E43').Value = '  -1.38%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.17'
$ws.Range('E44').Value = '  +5.61%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.32'
$ws.Range('E45').Value = '  -3.01%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0496'
$ws.Range('E46').Value = '  -4.17%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.959.48'
$ws.Range('E47').Value = '  -1.29%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.86'
$ws.Range('E48').Value = '  -3.71%  '

# Row 49
$ws.Range('E49').Value = '  +0.09%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '102.10'
$ws.Range('E50').Value = '  -3.43%  '

# Row 51
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '48.65'
$ws.Range('E51').Value = '  -3.08%  '

